$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 6 ("STGraph - Operations"): split the isValid(...) formula paragraph
# into three runs, the middle one italic.
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$tr6 = $s6.Shapes.Item(2).TextFrame.TextRange

# NB: the read-back of TextRange.Text mangles the U+2026 ellipsis character
# in this host, so anchor the search on an ellipsis-free substring and use
# the (known, fixed) character length of the original run instead of
# searching for the full string.
$oldFormulaAnchor = "isValid(Path(ni, "
$oldFormulaLength = 67
$part1 = "isValid(Path(ni, …, nk))  ⇔  ("
$part2 = "{j=i}^{k-1} I"
$part3 = "{e(n_j, n_{j+1})} ), where (I_e = [t_a,, t_b[)"
$combined = $part1 + $part2 + $part3

$full6 = $tr6.Text
$formulaStart = $full6.IndexOf($oldFormulaAnchor)

$formulaRange = $tr6.Characters($formulaStart + 1, $oldFormulaLength)
$formulaRange.Text = $combined

$italicStart = $formulaStart + 1 + $part1.Length
$italicRange = $tr6.Characters($italicStart, $part2.Length)
$italicRange.Font.Italic = $true

# ---------------------------------------------------------------------------
# Slide 7 ("Limitations and Future works"): drop the stray space before the
# trailing semicolons and append a new bullet.
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$tr7 = $s7.Shapes.Item(2).TextFrame.TextRange

$old7a = "No support for cross time-series operations ;"
$new7a = "No support for cross time-series operations;"
$full7 = $tr7.Text
$idx7a = $full7.IndexOf($old7a)
$range7a = $tr7.Characters($idx7a + 1, $old7a.Length)
$range7a.Text = $new7a

$old7b = "Query to AsterixDB should be asynchronous ;"
$new7b = "Query to AsterixDB should be asynchronous;"
$full7b = $tr7.Text
$idx7b = $full7b.IndexOf($old7b)
$range7b = $tr7.Characters($idx7b + 1, $old7b.Length)
$range7b.Text = $new7b

$null = $tr7.InsertAfter("`r" + "AsterixDB implementation could be far optimized and its full capabilities integrated;")
